$wb = $excel.ActiveWorkbook

# Update the "Metadata" sheet: URL, Version, Date, Publisher
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-received-date"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# Update the "Elements" sheet: clear the Constraint(s) value for the root "Extension" row (row 2, column AI)
$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("AI2").Value = ""

# The Extension.url row's Fixed Value mirrors the StructureDefinition URL, so it
# must be updated to the new URL too (it shares the same underlying string).
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-received-date"
